$wb = $excel.ActiveWorkbook

# ---- Status sheet ----
$status = $wb.Worksheets.Item("Status")

# Header row keeps "Deliverable" / "Date" / "Status (%)" as before; only
# the Date column values change from real dates to text-like date labels,
# and column B becomes wider to fit them.
$status.Range("B2").Value = "1st April 2023"
$status.Range("B3").Value = "1st May 2023"
$status.Range("B4").Value = "1st July 2023"
$status.Range("B5").Value = "1st August 2023"
$status.Range("B6").Value = "30th September 2023"

$status.Columns.Item(2).ColumnWidth = 22.4609375

$status.Range("B7").Select()

# ---- Tasks sheet ----
$tasks = $wb.Worksheets.Item("Tasks")

$tasks.Range("B1").Value = "Status (%)"

# Drop the percent number format and store plain numbers instead of
# fractions (0.5 -> 50, 0.4 -> 40, etc.)
$tasks.Range("B2:B6").NumberFormat = "General"
$tasks.Range("B2").Value = 50
$tasks.Range("B3").Value = 50
$tasks.Range("B4").Value = 40
$tasks.Range("B5").Value = 0
$tasks.Range("B6").Value = 0

$tasks.Range("H3").Select()
